$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 538.7143
$ws.Range("I2").Value = 354.2
$ws.Range("K2").Value = 354.2
$ws.Range("M2").Value = -241.2
$ws.Range("H18").Value = 724.73334
$ws.Range("I18").Value = 561
$ws.Range("J18").Value = 1175
$ws.Range("K18").Value = 561
$ws.Range("L18").Value = 1175
$ws.Range("M18").Value = -277
$ws.Range("N18").Value = -1743
$ws.Range("H21").Value = 65106.43
$ws.Range("I21").Value = 64149
$ws.Range("J21").Value = 67500
$ws.Range("K21").Value = 64149
$ws.Range("L21").Value = 67500
$ws.Range("M21").Value = -63681
$ws.Range("N21").Value = -68436
$ws.Range("H23").Value = 65106.43
$ws.Range("I23").Value = 64149
$ws.Range("J23").Value = 67500
$ws.Range("K23").Value = 64149
$ws.Range("L23").Value = 67500
$ws.Range("M23").Value = -63915
$ws.Range("N23").Value = -67968
$ws.Range("H34").Value = 3969.6155
$ws.Range("I34").Value = 2127.7273
$ws.Range("J34").Value = 14100
$ws.Range("K34").Value = 2127.7273
$ws.Range("L34").Value = 14100
$ws.Range("M34").Value = -1924.7273
$ws.Range("N34").Value = -14506
$ws.Range("H36").Value = 3969.6155
$ws.Range("I36").Value = 2127.7273
$ws.Range("J36").Value = 14100
$ws.Range("K36").Value = 2127.7273
$ws.Range("L36").Value = 14100
$ws.Range("M36").Value = -1412.7273
$ws.Range("N36").Value = -15530
$ws.Range("H70").Value = 2500.4707
$ws.Range("I70").Value = 2128.8572
$ws.Range("J70").Value = 2760.6
$ws.Range("K70").Value = 6386.571599999999
$ws.Range("L70").Value = 8281.799999999999
$ws.Range("M70").Value = -6116.571599999999
$ws.Range("N70").Value = -8821.799999999999
$ws.Range("H73").Value = 2500.4707
$ws.Range("I73").Value = 2128.8572
$ws.Range("J73").Value = 2760.6
$ws.Range("K73").Value = 6386.571599999999
$ws.Range("L73").Value = 8281.799999999999
$ws.Range("M73").Value = -5450.571599999999
$ws.Range("N73").Value = -10153.8
$ws.Range("H76").Value = 34217.938
$ws.Range("I76").Value = 61987.53
$ws.Range("J76").Value = 4712.75
$ws.Range("K76").Value = 61987.53
$ws.Range("L76").Value = 4712.75
$ws.Range("M76").Value = -61672.53
$ws.Range("N76").Value = -5342.75
$ws.Range("H79").Value = 34217.938
$ws.Range("I79").Value = 61987.53
$ws.Range("J79").Value = 4712.75
$ws.Range("K79").Value = 61987.53
$ws.Range("L79").Value = 4712.75
$ws.Range("M79").Value = -60895.53
$ws.Range("N79").Value = -6896.75
$ws.Range("H98").Value = 1703.1111
$ws.Range("I98").Value = 1287
$ws.Range("J98").Value = 2535.3333
$ws.Range("K98").Value = 1287
$ws.Range("L98").Value = 2535.3333
$ws.Range("M98").Value = 211
$ws.Range("N98").Value = -5531.3333
$ws.Range("H122").Value = 1703.1111
$ws.Range("I122").Value = 1287
$ws.Range("J122").Value = 2535.3333
$ws.Range("K122").Value = 3861
$ws.Range("L122").Value = 7605.999899999999
$ws.Range("M122").Value = -1411
$ws.Range("N122").Value = -12505.9999
$ws.Range("H129").Value = 1032
$ws.Range("I129").Value = 816
$ws.Range("J129").Value = 1045.3334
$ws.Range("K129").Value = 2448
$ws.Range("L129").Value = 3136.0002
$ws.Range("M129").Value = 2552
$ws.Range("N129").Value = -13136.0002
$ws.Range("H132").Value = 1478.3513
$ws.Range("I132").Value = 1356.7941
$ws.Range("J132").Value = 2856
$ws.Range("K132").Value = 4070.3823
$ws.Range("L132").Value = 8568
$ws.Range("M132").Value = -1540.3823
$ws.Range("N132").Value = -13628

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4361.787
$ws.Range("I61").Value = 5720.346
$ws.Range("K61").Value = 5720.346
$ws.Range("M61").Value = -5508.346
$ws.Range("H132").Value = 2439.15
$ws.Range("I132").Value = 1278.6578
$ws.Range("J132").Value = 4443.636
$ws.Range("K132").Value = 3835.9734
$ws.Range("L132").Value = 13330.908
$ws.Range("M132").Value = -1305.9734
$ws.Range("N132").Value = -18390.908
$ws.Range("H136").Value = 4361.787
$ws.Range("I136").Value = 5720.346
$ws.Range("K136").Value = 17161.038
$ws.Range("M136").Value = -14611.038

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 4638.25
$ws.Range("I29").Value = 3777.6667
$ws.Range("K29").Value = 3777.6667
$ws.Range("M29").Value = -3488.6667
$ws.Range("H134").Value = 4656.55
$ws.Range("I134").Value = 5877
$ws.Range("J134").Value = 2825.875
$ws.Range("K134").Value = 17631
$ws.Range("L134").Value = 8477.625
$ws.Range("M134").Value = -15096
$ws.Range("N134").Value = -13547.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2264.5715
$ws.Range("I16").Value = 1944.4286
$ws.Range("J16").Value = 2584.7144
$ws.Range("K16").Value = 1944.4286
$ws.Range("L16").Value = 2584.7144
$ws.Range("M16").Value = -1657.4286
$ws.Range("N16").Value = -3158.7144
$ws.Range("H113").Value = 2264.5715
$ws.Range("I113").Value = 1944.4286
$ws.Range("J113").Value = 2584.7144
$ws.Range("K113").Value = 1944.4286
$ws.Range("L113").Value = 2584.7144
$ws.Range("M113").Value = 225.5714
$ws.Range("N113").Value = -6924.7144
$ws.Range("H134").Value = 3102.8647
$ws.Range("I134").Value = 3305.6296
$ws.Range("J134").Value = 2555.4
$ws.Range("K134").Value = 9916.888800000001
$ws.Range("L134").Value = 7666.200000000001
$ws.Range("M134").Value = -7381.888800000001
$ws.Range("N134").Value = -12736.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 5658
$ws.Range("I139").Value = 14003.75
$ws.Range("J139").Value = 3185.1853
$ws.Range("K139").Value = 42011.25
$ws.Range("L139").Value = 9555.555899999999
$ws.Range("M139").Value = -36871.25
$ws.Range("N139").Value = -19835.5559

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 262835
$ws.Range("I12").Value = 324812.5
$ws.Range("J12").Value = 14925
$ws.Range("K12").Value = 324812.5
$ws.Range("L12").Value = 14925
$ws.Range("M12").Value = -324672.5
$ws.Range("N12").Value = -15205
$ws.Range("H14").Value = 2006600.9
$ws.Range("I14").Value = 6666668
$ws.Range("J14").Value = 9429.286
$ws.Range("K14").Value = 6666668
$ws.Range("L14").Value = 9429.286
$ws.Range("M14").Value = -6666500
$ws.Range("N14").Value = -9765.286
$ws.Range("H80").Value = 2510.6
$ws.Range("I80").Value = 2510.3125
$ws.Range("J80").Value = 2511.111
$ws.Range("K80").Value = 2510.3125
$ws.Range("L80").Value = 2511.111
$ws.Range("M80").Value = -1512.3125
$ws.Range("N80").Value = -4507.111
$ws.Range("H83").Value = 2510.6
$ws.Range("I83").Value = 2510.3125
$ws.Range("J83").Value = 2511.111
$ws.Range("K83").Value = 12551.5625
$ws.Range("L83").Value = 12555.555
$ws.Range("M83").Value = -7559.5625
$ws.Range("N83").Value = -22539.555
$ws.Range("H97").Value = 1142.3793
$ws.Range("I97").Value = 1147.4642
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 1147.4642
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -651.4641999999999
$ws.Range("N97").Value = -1992
$ws.Range("H113").Value = 125002270
$ws.Range("I113").Value = 200001230
$ws.Range("J113").Value = 4004.3333
$ws.Range("K113").Value = 200001230
$ws.Range("L113").Value = 4004.3333
$ws.Range("M113").Value = -199999060
$ws.Range("N113").Value = -8344.3333
$ws.Range("H122").Value = 16342899
$ws.Range("I122").Value = 2316798
$ws.Range("K122").Value = 6950394
$ws.Range("M122").Value = -6947944
$ws.Range("H132").Value = 24318.088
$ws.Range("I132").Value = 41015.08
$ws.Range("J132").Value = 2612
$ws.Range("K132").Value = 123045.24
$ws.Range("L132").Value = 7836
$ws.Range("M132").Value = -120515.24
$ws.Range("N132").Value = -12896

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 49579.094
$ws.Range("I7").Value = 57521.39
$ws.Range("J7").Value = 1925.3334
$ws.Range("K7").Value = 57521.39
$ws.Range("L7").Value = 1925.3334
$ws.Range("M7").Value = -57409.39
$ws.Range("N7").Value = -2149.3334
$ws.Range("H40").Value = 34484828
$ws.Range("I40").Value = 45456690
$ws.Range("J40").Value = 1829.2858
$ws.Range("K40").Value = 45456690
$ws.Range("L40").Value = 1829.2858
$ws.Range("M40").Value = -45456554
$ws.Range("N40").Value = -2101.2858
$ws.Range("H61").Value = 1679.091
$ws.Range("I61").Value = 1607.7778
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1607.7778
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1405.7778
$ws.Range("N61").Value = -2404
$ws.Range("H113").Value = 1679.091
$ws.Range("I113").Value = 1607.7778
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1607.7778
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 562.2221999999999
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 3019467.5
$ws.Range("I122").Value = 3574931
$ws.Range("J122").Value = 1432428.6
$ws.Range("K122").Value = 10724793
$ws.Range("L122").Value = 4297285.800000001
$ws.Range("M122").Value = -10722343
$ws.Range("N122").Value = -4302185.800000001
$ws.Range("H126").Value = 49579.094
$ws.Range("I126").Value = 57521.39
$ws.Range("J126").Value = 1925.3334
$ws.Range("K126").Value = 172564.17
$ws.Range("L126").Value = 5776.0002
$ws.Range("M126").Value = -170094.17
$ws.Range("N126").Value = -10716.0002
$ws.Range("H132").Value = 9528936
$ws.Range("I132").Value = 14499413
$ws.Range("J132").Value = 2189.0833
$ws.Range("K132").Value = 43498239
$ws.Range("L132").Value = 6567.249899999999
$ws.Range("M132").Value = -43495709
$ws.Range("N132").Value = -11627.2499
$ws.Range("H136").Value = 5255.864
$ws.Range("I136").Value = 4133.0933
$ws.Range("K136").Value = 12399.2799
$ws.Range("M136").Value = -9849.279900000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1559.9736
$ws.Range("I113").Value = 1252.3704
$ws.Range("J113").Value = 2315
$ws.Range("K113").Value = 3757.1112
$ws.Range("L113").Value = 6945
$ws.Range("M113").Value = -1587.1112
$ws.Range("N113").Value = -11285
$ws.Range("H122").Value = 2083.3333
$ws.Range("I122").Value = 2083.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6249.999899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3799.999899999999
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 1325.8334
$ws.Range("I126").Value = 975
$ws.Range("K126").Value = 2925
$ws.Range("M126").Value = -455
$ws.Range("H132").Value = 2231.2
$ws.Range("I132").Value = 1739.35
$ws.Range("J132").Value = 4198.6
$ws.Range("K132").Value = 5218.049999999999
$ws.Range("L132").Value = 12595.8
$ws.Range("M132").Value = -2688.049999999999
$ws.Range("N132").Value = -17655.8
$ws.Range("H136").Value = 2109.8545
$ws.Range("I136").Value = 2374.6775
$ws.Range("J136").Value = 1767.7916
$ws.Range("K136").Value = 7124.032499999999
$ws.Range("L136").Value = 5303.3748
$ws.Range("M136").Value = -4574.032499999999
$ws.Range("N136").Value = -10403.3748
